# BE: login, logout, email+code
# Adds 5 new test rows (12-16) to the "REPORT TESTS BACK-END" sheet:
#   12: view "email.blade" / vista che consente di impostare il corpo della mail
#   13-14 (merged A/B): recupero code nella view / inserimento link nella vista
#         con code per reset password  -> failed ("inserire il route nel link"),
#         then completato ("soluzione -> {{route...}}")
#   15: metodo: login / metodo per autenticazione  -> completato / soluzione-> Auth
#   16: metodo: logout / metodo per logout (cancellazione sessione creata nel login)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122
# xlCenter (used for VerticalAlignment / HorizontalAlignment)
$xlCenter = -4108

# --- Row 12: plain "completato" row, same look as rows 3-6 -----------------
$ws.Range("A3:E3").Copy()
$ws.Range("A12:E12").PasteSpecial($xlPasteFormats)

$ws.Range("A12").Value = "view ""email.blade"""
$ws.Range("B12").Value = "vista che consente di impostare il corpo della mail"
$ws.Range("C12").Value = 44537
$ws.Range("D12").Value = "completato"

# --- Rows 13-14: merged A/B block, mirrors rows 7-8 (failed -> completato) -
$ws.Range("A7:E7").Copy()
$ws.Range("A13:E13").PasteSpecial($xlPasteFormats)
$ws.Range("A8:E8").Copy()
$ws.Range("A14:E14").PasteSpecial($xlPasteFormats)

$ws.Range("A13:A14").Merge()
$ws.Range("B13:B14").Merge()

$ws.Range("A13").Value = "recupero code nella view"
$ws.Range("B13").Value = "inserimento link nella vista con code per reset password"
$ws.Range("C13").Value = 44537
$ws.Range("D13").Value = "failed"
$ws.Range("E13").Value = "inserire il route nel link"

$ws.Range("C14").Value = 44537
$ws.Range("D14").Value = "completato"
$ws.Range("E14").Value = "soluzione -> {{route...}}"

# Dates in this block use dd/mm/yy (same as rows 3-6), not the d/m/yy used by
# the row 7/8 block they were copied from.
$ws.Range("C13").NumberFormat = "dd/mm/yy"
$ws.Range("C14").NumberFormat = "dd/mm/yy"

# E14 now holds text, so it should not stay blank-styled - match the style
# used elsewhere for a populated Note cell on an unfilled row (e.g. E11).
$ws.Range("E11").Copy()
$ws.Range("E14").PasteSpecial($xlPasteFormats)
$ws.Range("E14").Value = "soluzione -> {{route...}}"

# --- Row 15: login ----------------------------------------------------------
$ws.Range("A3:E3").Copy()
$ws.Range("A15:E15").PasteSpecial($xlPasteFormats)

$ws.Range("A15").Value = "metodo: login"
$ws.Range("B15").Value = "metodo per autenticazione "
$ws.Range("C15").Value = 44537
$ws.Range("D15").Value = "completato"
$ws.Range("D15").VerticalAlignment = $xlCenter

$ws.Range("E11").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)
$ws.Range("E15").Value = "soluzione-> Auth"

# --- Row 16: logout ----------------------------------------------------------
$ws.Range("A3:E3").Copy()
$ws.Range("A16:E16").PasteSpecial($xlPasteFormats)

$ws.Range("A16").Value = "metodo: logout"
$ws.Range("B16").Value = "metodo per logout (cancellazione sessione creata nel login)"
$ws.Range("C16").Value = 44537
$ws.Range("D16").Value = "completato"
$ws.Range("D16").VerticalAlignment = $xlCenter
